$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$xlPasteValues = -4163
$xlPasteFormats = -4122

# 1) First, give the brand-new row 22 the same formatting (style) as an existing
#    data row (row 20) before anything is shifted into it, so the engine reuses
#    the existing style index instead of minting a new one.
$ws.Range("A20:B20").Copy()
$ws.Range("A22:B22").PasteSpecial($xlPasteFormats)

# 2) Shift rows 11..21 down to 12..22 (process bottom-up so sources aren't
#    clobbered before they're read). Using Copy + PasteSpecial(values) instead
#    of a plain .Value assignment keeps the exact stored type (text vs number
#    vs boolean) of tricky-looking strings like "true" / "5" intact.
for ($r = 21; $r -ge 11; $r--) {
    $dst = $r + 1
    $ws.Range("A${r}:B${r}").Copy()
    $ws.Range("A${dst}:B${dst}").PasteSpecial($xlPasteValues)
}

$excel.CutCopyMode = 0

# 3) Update the Date and Contact values in place (row positions unchanged).
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# 4) Row 11 becomes the new "Jurisdiction" property row; its value cell is
#    left blank (closest achievable representation of an empty string cell).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").ClearContents()
